{"js": "const pairs = [\n  [\"2025-02-01 Saturday\", \"2025-02-02 Sunday\"],\n  [\"575\u00f79=\", \"490\u00f74=\"],\n  [\"619\u00f72=\", \"191\u00f76=\"],\n  [\"746\u00f72=\", \"769\u00f78=\"],\n  [\"241\u00f75=\", \"183\u00f76=\"],\n  [\"615\u00f77=\", \"944\u00f78=\"],\n  [\"616\u00f73=\", \"295\u00f72=\"],\n  [\"398\u00f78=\", \"905\u00f74=\"],\n  [\"910\u00f73=\", \"496\u00f74=\"],\n  [\"216\u00f79=\", \"242\u00f74=\"],\n  [\"168\u00f76=\", \"741\u00f76=\"],\n  [\"196\u00f76=\", \"276\u00f74=\"],\n  [\"409\u00f78=\", \"525\u00f78=\"],\n  [\"569\u00f73=\", \"254\u00f79=\"],\n  [\"461\u00f73=\", \"997\u00f75=\"],\n  [\"659\u00f79=\", \"147\u00f78=\"],\n  [\"235\u00f72=\", \"460\u00f78=\"],\n  [\"313\u00f76=\", \"821\u00f75=\"],\n  [\"408\u00f72=\", \"839\u00f73=\"],\n  [\"847\u00f75=\", \"216\u00f73=\"],\n  [\"591\u00f75=\", \"647\u00f73=\"],\n  [\"710\u00f77=\", \"707\u00f77=\"],\n  [\"927\u00f77=\", \"310\u00f73=\"],\n  [\"865\u00f79=\", \"341\u00f78=\"],\n  [\"184\u00f79=\", \"113\u00f78=\"],\n  [\"122\u00f78=\", \"269\u00f72=\"],\n];\n\nconst body = context.document.body;\nfor (const [before, after] of pairs) {\n  const results = body.search(before, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(after, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n  ,@(\"2025-02-01 Saturday\", \"2025-02-02 Sunday\")\n  ,@(\"575\u00f79=\", \"490\u00f74=\")\n  ,@(\"619\u00f72=\", \"191\u00f76=\")\n  ,@(\"746\u00f72=\", \"769\u00f78=\")\n  ,@(\"241\u00f75=\", \"183\u00f76=\")\n  ,@(\"615\u00f77=\", \"944\u00f78=\")\n  ,@(\"616\u00f73=\", \"295\u00f72=\")\n  ,@(\"398\u00f78=\", \"905\u00f74=\")\n  ,@(\"910\u00f73=\", \"496\u00f74=\")\n  ,@(\"216\u00f79=\", \"242\u00f74=\")\n  ,@(\"168\u00f76=\", \"741\u00f76=\")\n  ,@(\"196\u00f76=\", \"276\u00f74=\")\n  ,@(\"409\u00f78=\", \"525\u00f78=\")\n  ,@(\"569\u00f73=\", \"254\u00f79=\")\n  ,@(\"461\u00f73=\", \"997\u00f75=\")\n  ,@(\"659\u00f79=\", \"147\u00f78=\")\n  ,@(\"235\u00f72=\", \"460\u00f78=\")\n  ,@(\"313\u00f76=\", \"821\u00f75=\")\n  ,@(\"408\u00f72=\", \"839\u00f73=\")\n  ,@(\"847\u00f75=\", \"216\u00f73=\")\n  ,@(\"591\u00f75=\", \"647\u00f73=\")\n  ,@(\"710\u00f77=\", \"707\u00f77=\")\n  ,@(\"927\u00f77=\", \"310\u00f73=\")\n  ,@(\"865\u00f79=\", \"341\u00f78=\")\n  ,@(\"184\u00f79=\", \"113\u00f78=\")\n  ,@(\"122\u00f78=\", \"269\u00f72=\")\n)\n\nforeach ($pair in $pairs) {\n  $findText = $pair[0]\n  $replaceText = $pair[1]\n  $range = $d.Content\n  $find = $range.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}"}
